$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "1. a) Describe the characteristics and benefits of cloud computing. Explain how it differs from traditional computing models. [7]" "1. a) Describe the characteristics and benefits of cloud computing. Explain how cloud computing can be applied in the healthcare industry. [7]"

Replace-Text "b) Discuss the NIST definition of cloud computing and its reference model, highlighting the key components and their interactions. [8]" "b) Compare and contrast the different cloud service models (IaaS, PaaS, SaaS) and deployment models (public, private, hybrid, community). [8]"

Replace-Text "2. a) Explain the concept of virtualization in cloud computing, including its types and benefits. [7]" "2. a) Design an architecture for a static web hosting system on a Linux instance in a cloud environment. [7]"

Replace-Text "b) Design a virtual private cloud (VPC) architecture for a small business, including subnets, security groups, and network ACLs. [8]" "b) Explain the concept of virtualization and its role in cloud computing. Discuss the differences between Type 1 and Type 2 hypervisors. [8]"

Replace-Text "3. a) Compare and contrast the different cloud service models (IaaS, PaaS, SaaS) and deployment models (public, private, hybrid, community). [7]" "3. a) What are the key components of a cloud storage system? Explain the differences between object storage, file storage, and block storage. [7]"

Replace-Text "b) A company has 1000 users and wants to deploy a cloud-based email service. Calculate the total cost of ownership (TCO) for a 3-year period, assuming an annual subscription fee of `$10 per user for SaaS, `$50 per user for PaaS, and `$100 per user for IaaS. [8]" "b) A company has 1000 employees and wants to implement a cloud-based database system. The estimated monthly data storage requirement is 500 GB. Calculate the total cost of using a cloud storage service that charges `$0.10 per GB per month. [8]"

Replace-Text "4. a) Discuss the concept of auto-scaling in cloud computing, including its benefits and challenges. [7]" "4. a) Explain the concept of auto-scaling in cloud computing. Discuss the different scaling strategies and their applications. [7]"

Replace-Text "b) A web application experiences a sudden surge in traffic, with the number of requests increasing from 100 to 1000 per minute. Design an auto-scaling strategy to handle this increase, including the use of load balancers and scaling rules. [8]" "b) Design an architecture to implement auto-scaling and load balancing to distribute the workload of a dynamic web application in a cloud environment. [8]"

Replace-Text "5. a) Explain the concept of cloud storage models, including block storage, file storage, and object storage. [7]" "5. a) What are the security threats and challenges in cloud computing? Explain the concept of security as a service (SaaS) and its importance in cloud security. [7]"

Replace-Text "b) A company has 10 TB of data stored in a cloud-based object storage service, with an average retrieval rate of 1000 requests per hour. Calculate the total cost of storage and retrieval for a 1-year period, assuming a storage cost of `$0.01 per GB per month and a retrieval cost of `$0.005 per 1000 requests. [8]" "b) A company wants to implement a cloud-based identity and access management (IAM) system. Explain the key components and benefits of an IAM system. [8]"

Replace-Text "6. a) Discuss the concept of security in cloud computing, including the types of threats and vulnerabilities. [7]" "6. a) Explain the concept of cloud monitoring and its importance in cloud computing. Discuss the different types of cloud monitoring tools and their applications. [7]"

Replace-Text "b) Design a security architecture for a cloud-based e-commerce application, including the use of firewalls, intrusion detection systems, and encryption. [8]" "b) A company has a cloud-based application that requires a minimum of 5 instances to run during peak hours. The estimated monthly cost of each instance is `$100. Calculate the total monthly cost of using a cloud service that charges based on the number of instances used. [8]"

Replace-Text "i) Cloud security mechanisms" "i) Cloud Reference Model"

Replace-Text "ii) Cloud monitoring and management tools" "ii) Virtual Private Cloud (VPC)"

Replace-Text "iii) Cloud deployment models [2 x 5] [10]" "iii) Service Level Agreement (SLA) [2 x 5]"
